$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.39%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.29%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.601"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.25%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05891"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.08%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-0.78%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8522"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.13%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9420"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.15%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.49%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05049"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'40.88%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07077"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03104"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.13%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09127"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.31%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001537"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.51%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006071"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.006076"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.33%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.493"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.183"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.26%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'0.3055"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.85%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1270"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.78%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.938"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'11.63%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04268"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.47%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.35%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004287"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-5.22%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'29.77%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.44%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006295"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'6.90%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1099"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.28%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.10%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01415"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'31.28%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005366"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.20%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.10%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05100"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-53.24%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'11,563.61%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
